$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: old "A lot of work..." paragraph becomes a new underlined
# "Problem" heading (+ blank line) followed by the reworded paragraph.
# ---------------------------------------------------------------
$needle1 = "A lot of work has been done in identifying sentiment in social media posts based on text analysis.   However, much social media communication is image-based, and classifying sentiment from these images has proved to be a larger challenge, in part due to the difficulty in obtaining sufficient training data."

$replace1 = "Problem^p^p" + `
    "Much work has been done using text to identify sentiment in social media posts.   " + `
    "However, much social media communication is image-based.  " + `
    "Classifying sentiment from these images has proved to be a larger challenge, " + `
    "partially because of the difficulty in obtaining sufficient training data."

$found1 = $d.Content.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)
Write-Output ("Edit1 Found=" + $found1)

# ---------------------------------------------------------------
# Edit 2: the old "I propose..." paragraph, the blank line after it,
# and the old "Using the categories..." paragraph (which carried the
# footnote reference) are all replaced by two new underlined headings
# ("Approach and Data", "Methodology") plus three reworded paragraphs.
# This also removes the footnote, since its reference run sits inside
# the replaced span.
# ---------------------------------------------------------------
$needle2 = "I propose to develop an unsupervised image sentiment classifier based on crawled Twitter data.  In order to classify image sentiment, I will generate a sentiment label based on text-based sentiment analysis, taking into account social/graph-based input (i.e., the text analysis will not be solely based on the original poster of an image, but also on retweets, responses, etc.).`r`rUsing the categories generated above, I will implement a Neural-Network based predictor on the images retrieved from my Twitter crawl.   The final predictor will be tested against images that have been sentiment-scored by crowd sourcing."

$replace2 = "Approach and Data^p^p" + `
    "I propose to develop an unsupervised image sentiment classifier based on streamed Twitter data.  ^p^p" + `
    "Methodology^p^p" + `
    "First, I will generate a sentiment label based on text-based sentiment analysis, taking into account both the original tweet and responses.^p^p" + `
    "I will use these text-derived sentiments to train a Neural-Network based predictor on the tweeted images.  "

$found2 = $d.Content.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)
Write-Output ("Edit2 Found=" + $found2)

# ---------------------------------------------------------------
# Underline the three new section headings. We re-fetch the owning
# paragraph via the document's Paragraphs collection (by numeric
# Index) so the whole paragraph range -- including its paragraph
# mark -- picks up the underline (matching pPr/rPr + run rPr).
# ---------------------------------------------------------------
foreach ($heading in @("Problem", "Approach and Data", "Methodology")) {
    $hr = $d.Content
    $hf = $hr.Find
    $hf.ClearFormatting()
    $ok = $hf.Execute($heading, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $idx = $hr.Paragraphs.Item(1).Index
        $d.Paragraphs.Item($idx).Range.Font.Underline = 1
    }
    Write-Output ("Underline '" + $heading + "' Found=" + $ok)
}

# ---------------------------------------------------------------
# Re-create the _GoBack bookmark at its new location: right after
# "text-derived" and before " sentiments to train" in the last
# paragraph of the Methodology section.
# ---------------------------------------------------------------
$br = $d.Content
$bf = $br.Find
$bf.ClearFormatting()
$bok = $bf.Execute("I will use these text-derived", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Bookmark anchor Found=" + $bok)
if ($bok) {
    $br.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $br) | Out-Null
}

Write-Output ("FootnotesCount=" + $d.Footnotes.Count)
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ($i.ToString() + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}
